$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "25.50"); force text
# formatting on just the rows being rewritten so Excel keeps the new values as
# literal strings instead of coercing them to numbers / normalizing digits.
$ws.Range("D2:D10").NumberFormat = "@"
$ws.Range("D12:D17").NumberFormat = "@"
$ws.Range("D19:D21").NumberFormat = "@"
$ws.Range("D23:D36").NumberFormat = "@"
$ws.Range("D38:D40").NumberFormat = "@"
$ws.Range("D42:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.439.89"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "2.128.99"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "348.24"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "0.5226"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("D8").Value = "0.4493"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "54.13"
$ws.Range("E9").Value = "  +3.99%  "
$ws.Range("D10").Value = "0.09399"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "25.55"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "8.720"
$ws.Range("E13").Value = "  +8.61%  "
$ws.Range("D14").Value = "2.132.31"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "6.985"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").Value = "103.31"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "0.00001176"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "21.65"
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").Value = "0.06718"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "6.346"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "30.438.63"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "12.81"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "2.339"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "2.373.19"
$ws.Range("D27").Value = "22.31"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "2.562"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "163.59"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "134.75"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "1.170"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "1.812"
$ws.Range("E32").Value = "  +12.19%  "
$ws.Range("D33").Value = "0.1061"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "6.908"
$ws.Range("E34").Value = "  +12.68%  "
$ws.Range("D35").Value = "6.324"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").Value = "3.958"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +6.17%  "
$ws.Range("D38").Value = "0.02649"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "0.06885"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "0.7168"
$ws.Range("E40").Value = "  +4.33%  "
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("D42").Value = "0.2261"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").Value = "1.342"
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("D44").Value = "0.6980"
$ws.Range("E44").Value = "  +5.45%  "
$ws.Range("D45").Value = "14.85"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("D46").Value = "2.402"
$ws.Range("E46").Value = "  +5.49%  "
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "1.276"
$ws.Range("E48").Value = "  +9.10%  "
$ws.Range("D49").Value = "3.640"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "0.00000000349"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "84.10"
$ws.Range("E51").Value = "  +3.09%  "
